# Developed Change Email System - allow users to change their email during login.
# "Account System: Change Email" moves from "Not Started" (column A) to "Done"
# (column C). "User Idle System" is dropped from "Doing" (column B). A new
# "Account System: Activity Log" task is added at the bottom of "Not Started".

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item(1)

# Remove the finished "Account System: Change Email" card from the "Not
# Started" column (A6) and shift the remaining cards (A7:A10) up one row,
# then add the new task at the end of the list.
$a7 = $ws.Range("A7").Value()
$a8 = $ws.Range("A8").Value()
$a9 = $ws.Range("A9").Value()
$a10 = $ws.Range("A10").Value()
$ws.Range("A6").Value = $a7
$ws.Range("A7").Value = $a8
$ws.Range("A8").Value = $a9
$ws.Range("A9").Value = $a10
$ws.Range("A10").Value = "Account System: Activity Log"

# "User Idle System" is no longer tracked in "Doing".
$ws.Range("B4").ClearContents()

# Record "Account System: Change Email" as completed in the "Done" column.
$ws.Range("C14").Value = "Account System: Change Email"

# Match the author's last active selection.
$ws.Range("B10").Select()

$wb.Save()
